# Apply cryptos list data refresh as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.852.37"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "1.812.88"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'310.01"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "'0.4635"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("D8").Value = "'0.3693"
$ws.Range("E8").Value = "  -2.03%  "
$ws.Range("D9").Value = "'0.07347"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").Value = "'0.8697"
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").Value = "'20.40"
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("D12").Value = "1.832.68"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").Value = "'5.341"
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("D14").Value = "'0.07076"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "'6.509"
$ws.Range("E15").Value = "  -2.66%  "
$ws.Range("D16").Value = "'91.70"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "'0.000008706"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'14.67"
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").Value = "26.889.48"
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("D22").Value = "'5.340"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("E23").Value = "  -3.54%  "
$ws.Range("D24").Value = "2.113.22"
$ws.Range("E24").Value = "  +2.96%  "
$ws.Range("D25").Value = "'1.897"
$ws.Range("E25").Value = "  -2.05%  "
$ws.Range("D26").Value = "'151.88"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").Value = "'2.122"
$ws.Range("E28").Value = "  -5.84%  "
$ws.Range("D29").Value = "'5.295"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'115.24"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("D31").Value = "'0.08905"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "'0.7554"
$ws.Range("E32").Value = "  -3.36%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'2.936"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.149"
$ws.Range("E34").Value = "  -2.87%  "
$ws.Range("D35").Value = "'4.455"
$ws.Range("E35").Value = "  -1.67%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").Value = "'1.095"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("D39").Value = "'0.05252"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").Value = "'2.917"
$ws.Range("E40").Value = "  +0.75%  "
$ws.Range("D41").Value = "'0.5328"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").Value = "'7.192"
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("D43").Value = "'2.355"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("D44").Value = "'0.1660"
$ws.Range("E44").Value = "  -1.82%  "
$ws.Range("D45").Value = "'8.419"
$ws.Range("E45").Value = "  -2.27%  "
$ws.Range("D46").Value = "'0.4930"
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("D47").Value = "'10.38"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").Value = "'1.671"
$ws.Range("D50").Value = "'102.83"
$ws.Range("E50").Value = "  -2.69%  "
$ws.Range("D51").Value = "'0.06267"
$ws.Range("E51").Value = "  -1.05%  "
